$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 650984
$ws.Range("E2").Value = 32135
$ws.Range("F2").Value = 32135
$ws.Range("G2").Value = 13781
$ws.Range("H2").Value = 5567
$ws.Range("I2").Value = 6261
$ws.Range("J2").Value = -694
$ws.Range("K2").Value = 852522
$ws.Range("L2").Value = 399608
$ws.Range("M2").Value = 452914
$ws.Range("N2").Value = 415874
$ws.Range("O2").Value = 37040
$ws.Range("P2").Value = 4824
$ws.Range("Q2").Value = 34121
$ws.Range("R2").Value = -37452
$ws.Range("S2").Value = 1351
$ws.Range("T2").Value = 35055
$ws.Range("U2").Value = -934
$ws.Range("V2").Value = 274750
$ws.Range("W2").Value = 4.94
$ws.Range("X2").Value = 0.86
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 0.66
$ws.Range("AA2").Value = 88.23
$ws.Range("AB2").Value = 8717.040000000001
$ws.Range("AC2").Value = 7181
$ws.Range("AD2").Value = 38.36
$ws.Range("AE2").Value = 519887
$ws.Range("AF2").Value = 0.53
$ws.Range("AG2").Value = 8000
$ws.Range("AH2").Value = 2.9
$ws.Range("AI2").Value = 102.14
$ws.Range("AJ2").Value = 87186835
$ws.Range("D3").Value = 581923
$ws.Range("E3").Value = 24100
$ws.Range("F3").Value = 24100
$ws.Range("G3").Value = 1808
$ws.Range("H3").Value = -962
$ws.Range("I3").Value = 1806
$ws.Range("J3").Value = -2768
$ws.Range("K3").Value = 804088
$ws.Range("L3").Value = 353385
$ws.Range("M3").Value = 450702
$ws.Range("N3").Value = 412354
$ws.Range("O3").Value = 38349
$ws.Range("P3").Value = 4824
$ws.Range("Q3").Value = 76018
$ws.Range("R3").Value = -45347
$ws.Range("S3").Value = -22416
$ws.Range("T3").Value = 25602
$ws.Range("U3").Value = 50416
$ws.Range("V3").Value = 252615
$ws.Range("W3").Value = 4.14
$ws.Range("X3").Value = -0.17
$ws.Range("Y3").Value = 0.44
$ws.Range("Z3").Value = -0.12
$ws.Range("AA3").Value = 78.41
$ws.Range("AB3").Value = 8682.51
$ws.Range("AC3").Value = 2072
$ws.Range("AD3").Value = 80.36
$ws.Range("AE3").Value = 515470
$ws.Range("AF3").Value = 0.32
$ws.Range("AG3").Value = 8000
$ws.Range("AH3").Value = 4.8
$ws.Range("AI3").Value = 354.26
$ws.Range("AJ3").Value = 87186835
$ws.Range("D4").Value = 530835
$ws.Range("E4").Value = 28443
$ws.Range("F4").Value = 28443
$ws.Range("G4").Value = 14329
$ws.Range("H4").Value = 10482
$ws.Range("I4").Value = 13633
$ws.Range("J4").Value = -3151
$ws.Range("K4").Value = 797630
$ws.Range("L4").Value = 339246
$ws.Range("M4").Value = 458384
$ws.Range("N4").Value = 423734
$ws.Range("O4").Value = 34650
$ws.Range("P4").Value = 4824
$ws.Range("Q4").Value = 52694
$ws.Range("R4").Value = -37546
$ws.Range("S4").Value = -39510
$ws.Range("T4").Value = 23241
$ws.Range("U4").Value = 29453
$ws.Range("V4").Value = 228194
$ws.Range("W4").Value = 5.36
$ws.Range("X4").Value = 1.98
$ws.Range("Y4").Value = 3.26
$ws.Range("Z4").Value = 1.31
$ws.Range("AA4").Value = 74.01000000000001
$ws.Range("AB4").Value = 8824.889999999999
$ws.Range("AC4").Value = 15637
$ws.Range("AD4").Value = 16.47
$ws.Range("AE4").Value = 529683
$ws.Range("AF4").Value = 0.49
$ws.Range("AG4").Value = 8000
$ws.Range("AH4").Value = 3.11
$ws.Range("AI4").Value = 46.94
$ws.Range("AJ4").Value = 87186835
$ws.Range("D5").Value = 606551
$ws.Range("E5").Value = 46218
$ws.Range("F5").Value = 46218
$ws.Range("G5").Value = 41797
$ws.Range("H5").Value = 29735
$ws.Range("I5").Value = 27901
$ws.Range("J5").Value = 1834
$ws.Range("K5").Value = 790250
$ws.Range("L5").Value = 315610
$ws.Range("M5").Value = 474640
$ws.Range("N5").Value = 437329
$ws.Range("O5").Value = 37311
$ws.Range("P5").Value = 4824
$ws.Range("Q5").Value = 56073
$ws.Range("R5").Value = -38179
$ws.Range("S5").Value = -15655
$ws.Range("T5").Value = 22876
$ws.Range("U5").Value = 33197
$ws.Range("V5").Value = 211567
$ws.Range("W5").Value = 7.62
$ws.Range("X5").Value = 4.9
$ws.Range("Y5").Value = 6.48
$ws.Range("Z5").Value = 3.74
$ws.Range("AA5").Value = 66.48999999999999
$ws.Range("AB5").Value = 9218.26
$ws.Range("AC5").Value = 32001
$ws.Range("AD5").Value = 10.39
$ws.Range("AE5").Value = 546664
$ws.Range("AF5").Value = 0.61
$ws.Range("AG5").Value = 8000
$ws.Range("AH5").Value = 2.41
$ws.Range("AI5").Value = 22.94
$ws.Range("AJ5").Value = 87186835
$ws.Range("D6").Value = 649778
$ws.Range("E6").Value = 55426
$ws.Range("F6").Value = 55426
$ws.Range("G6").Value = 35628
$ws.Range("H6").Value = 18921
$ws.Range("I6").Value = 16906
$ws.Range("K6").Value = 782483
$ws.Range("L6").Value = 314887
$ws.Range("M6").Value = 467596
$ws.Range("N6").Value = 433713
$ws.Range("P6").Value = 4824
$ws.Range("Q6").Value = 58697
$ws.Range("R6").Value = -26480
$ws.Range("S6").Value = -31950
$ws.Range("T6").Value = 21356
$ws.Range("U6").Value = 37342
$ws.Range("V6").Value = 203040
$ws.Range("W6").Value = 8.529999999999999
$ws.Range("X6").Value = 2.91
$ws.Range("Y6").Value = 3.88
$ws.Range("Z6").Value = 2.41
$ws.Range("AA6").Value = 67.34
$ws.Range("AB6").Value = 9458.18
$ws.Range("AC6").Value = 19391
$ws.Range("AD6").Value = 12.53
$ws.Range("AE6").Value = 542133
$ws.Range("AF6").Value = 0.45
$ws.Range("AG6").Value = 10000
$ws.Range("AH6").Value = 4.12
$ws.Range("AI6").Value = 47.32
$ws.Range("AJ6").Value = 87186835
$ws.Range("D7").Value = 647004
$ws.Range("E7").Value = 41372
$ws.Range("G7").Value = 35130
$ws.Range("H7").Value = 23623
$ws.Range("I7").Value = 21758
$ws.Range("K7").Value = 793863
$ws.Range("L7").Value = 311434
$ws.Range("M7").Value = 482428
$ws.Range("N7").Value = 447368
$ws.Range("P7").Value = 4821
$ws.Range("Q7").Value = 57929
$ws.Range("R7").Value = -28886
$ws.Range("S7").Value = -22423
$ws.Range("T7").Value = 27742
$ws.Range("U7").Value = 27309
$ws.Range("W7").Value = 6.39
$ws.Range("X7").Value = 3.65
$ws.Range("Y7").Value = 4.94
$ws.Range("Z7").Value = 3
$ws.Range("AA7").Value = 64.56
$ws.Range("AC7").Value = 24956
$ws.Range("AD7").Value = 9.640000000000001
$ws.Range("AE7").Value = 558403
$ws.Range("AF7").Value = 0.43
$ws.Range("AG7").Value = 10150
$ws.Range("AH7").Value = 4.22
$ws.Range("AI7").Value = 40.67
$ws.Range("D8").Value = 652187
$ws.Range("E8").Value = 41067
$ws.Range("G8").Value = 36917
$ws.Range("H8").Value = 25919
$ws.Range("I8").Value = 23796
$ws.Range("K8").Value = 805994
$ws.Range("L8").Value = 307516
$ws.Range("M8").Value = 498477
$ws.Range("N8").Value = 461997
$ws.Range("P8").Value = 4821
$ws.Range("Q8").Value = 59567
$ws.Range("R8").Value = -38318
$ws.Range("S8").Value = -15025
$ws.Range("T8").Value = 32506
$ws.Range("U8").Value = 21551
$ws.Range("W8").Value = 6.3
$ws.Range("X8").Value = 3.97
$ws.Range("Y8").Value = 5.24
$ws.Range("Z8").Value = 3.24
$ws.Range("AA8").Value = 61.69
$ws.Range("AC8").Value = 27293
$ws.Range("AD8").Value = 8.1
$ws.Range("AE8").Value = 576663
$ws.Range("AF8").Value = 0.38
$ws.Range("AG8").Value = 10139
$ws.Range("AH8").Value = 4.59
$ws.Range("AI8").Value = 37.15
$ws.Range("D9").Value = 666023
$ws.Range("E9").Value = 44149
$ws.Range("G9").Value = 40347
$ws.Range("H9").Value = 28337
$ws.Range("I9").Value = 26235
$ws.Range("K9").Value = 822147
$ws.Range("L9").Value = 304763
$ws.Range("M9").Value = 517385
$ws.Range("N9").Value = 479243
$ws.Range("P9").Value = 4821
$ws.Range("Q9").Value = 61418
$ws.Range("R9").Value = -38759
$ws.Range("S9").Value = -14896
$ws.Range("T9").Value = 31941
$ws.Range("U9").Value = 24484
$ws.Range("W9").Value = 6.63
$ws.Range("X9").Value = 4.25
$ws.Range("Y9").Value = 5.58
$ws.Range("Z9").Value = 3.48
$ws.Range("AA9").Value = 58.9
$ws.Range("AC9").Value = 30091
$ws.Range("AD9").Value = 7.34
$ws.Range("AE9").Value = 598189
$ws.Range("AF9").Value = 0.37
$ws.Range("AG9").Value = 10167
$ws.Range("AH9").Value = 4.6
$ws.Range("AI9").Value = 33.79
